# daily auto push: 2025-09-30 13:37 UTC
# Append the new daily-tracker row (2025/09/30, 火, 20, 150) at the bottom
# of the sheet, extending the used range from A1:D40 to A1:D41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 41

# Column A holds a date-formatted literal like "2025/09/30", stored as plain
# text (inlineStr) in the source file rather than an actual date value.
# Force text entry first so Excel does not auto-convert the string into a
# date serial number, then drop back to the sheet's default ("Normal")
# style so the new cell matches the rest of the column (no explicit style).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/09/30"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "火"
$ws.Cells.Item($newRow, 3).Value = 20
$ws.Cells.Item($newRow, 4).Value = 150
